# Updates the cryptos price/volume table (and two name/link swaps) to match
# the latest scraped snapshot, per the "Updated cryptos list" GitHub Action commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.064.52'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '2.010.87'
$ws.Range("E3").Value = '  -2.03%  '

$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.20'
$ws.Range("E5").Value = '  -1.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.604'
$ws.Range("E6").Value = '  -1.41%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.82'
$ws.Range("E8").Value = '  -3.98%  '

$ws.Range("E9").Value = '  -2.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0780'
$ws.Range("E10").Value = '  -3.57%  '

$ws.Range("E11").Value = '  -5.16%  '

$ws.Range("D12").Value = '2.309.90'
$ws.Range("E12").Value = '  -2.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.02'
$ws.Range("E13").Value = '  -4.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.78'
$ws.Range("E14").Value = '  -5.00%  '

$ws.Range("E15").Value = '  -2.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.736'
$ws.Range("E16").Value = '  -2.93%  '

$ws.Range("D17").Value = '2.009.10'
$ws.Range("E17").Value = '  -2.16%  '

$ws.Range("D18").Value = '37.002.24'
$ws.Range("E18").Value = '  -0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  +2.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.22'
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("D21").Value = '0.0₃0812'
$ws.Range("E21").Value = '  -2.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.59'
$ws.Range("E22").Value = '  -2.03%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("E24").Value = '  +2.46%  '

$ws.Range("E25").Value = '  -5.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.45'
$ws.Range("E26").Value = '  -2.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.97'
$ws.Range("E27").Value = '  -6.00%  '

$ws.Range("E28").Value = '  -3.80%  '

$ws.Range("E29").Value = '  -2.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.29'
$ws.Range("E30").Value = '  -8.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.117'
$ws.Range("E31").Value = '  -1.43%  '

$ws.Range("E32").Value = '  -2.08%  '

# Row 33/34 swap: Hedera now ranks above InternetComputer(DFINITY)
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0600'
$ws.Range("E33").Value = '  -2.51%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.48'
$ws.Range("E34").Value = '  -3.02%  '

$ws.Range("E35").Value = '  -4.09%  '

$ws.Range("E36").Value = '  +2.06%  '

$ws.Range("E37").Value = '  +0.92%  '

$ws.Range("E38").Value = '  -4.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.33'
$ws.Range("E39").Value = '  -1.51%  '

$ws.Range("D40").Value = '1.454.28'
$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("E41").Value = '  -4.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.75'

# Row 43/44 swap: Cronos now ranks above HuobiToken
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0905'
$ws.Range("E43").Value = '  -4.33%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.75'
$ws.Range("E44").Value = '  -4.88%  '

$ws.Range("E45").Value = '  -3.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.91'
$ws.Range("E46").Value = '  -6.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.996'
$ws.Range("E48").Value = '  -2.66%  '

$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("D50").Value = '2.197.55'
$ws.Range("E50").Value = '  -2.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.59'
$ws.Range("E51").Value = '  -4.20%  '
